$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.260.40"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "'2.038.08"
$ws.Range("E3").Value = "  -2.41%  "

$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").Value = "'228.28"
$ws.Range("E5").Value = "  -2.34%  "

$ws.Range("D6").Value = "'0.609"
$ws.Range("E6").Value = "  -2.53%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -4.41%  "

$ws.Range("E9").Value = "  -3.10%  "

$ws.Range("D10").Value = "'0.0805"
$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("D11").Value = "'0.106"
$ws.Range("E11").Value = "  -2.08%  "

$ws.Range("D12").Value = "'2.340.34"
$ws.Range("E12").Value = "  -2.39%  "

$ws.Range("D13").Value = "'14.45"
$ws.Range("E13").Value = "  -4.71%  "

$ws.Range("D14").Value = "'20.40"
$ws.Range("E14").Value = "  -4.62%  "

$ws.Range("E15").Value = "  -4.10%  "

$ws.Range("E16").Value = "  -2.40%  "

$ws.Range("D17").Value = "'2.033.99"
$ws.Range("E17").Value = "  -2.70%  "

$ws.Range("D18").Value = "'37.134.47"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("E19").Value = "  -3.38%  "

$ws.Range("D20").Value = "'69.40"
$ws.Range("E20").Value = "  -2.67%  "

$ws.Range("D21").Value = "'0.0₃0836"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "'224.86"
$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").Value = "'2.36"
$ws.Range("E24").Value = "  -1.27%  "

$ws.Range("D25").Value = "'2.25"
$ws.Range("E25").Value = "  -6.16%  "

$ws.Range("D26").Value = "'9.45"
$ws.Range("E26").Value = "  -4.02%  "

$ws.Range("D27").Value = "'167.62"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("D28").Value = "'0.127"
$ws.Range("E28").Value = "  -6.65%  "

$ws.Range("E29").Value = "  -1.65%  "

$ws.Range("D30").Value = "'18.84"
$ws.Range("E30").Value = "  -3.58%  "

$ws.Range("E31").Value = "  -3.26%  "

$ws.Range("E32").Value = "  -5.16%  "

$ws.Range("D33").Value = "'4.54"
$ws.Range("E33").Value = "  -2.99%  "

$ws.Range("D34").Value = "'0.0607"
$ws.Range("E34").Value = "  -4.08%  "

$ws.Range("E35").Value = "  -5.12%  "

$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("D38").Value = "'3.16"
$ws.Range("E38").Value = "  -7.09%  "

$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("E40").Value = "  -7.77%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'1.481.56"
$ws.Range("E41").Value = "  +1.75%  "

$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.89"
$ws.Range("E42").Value = "  -1.49%  "

$ws.Range("D43").Value = "'16.73"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("D44").Value = "'0.0935"
$ws.Range("E44").Value = "  -4.23%  "

$ws.Range("D45").Value = "'95.08"
$ws.Range("E45").Value = "  -7.47%  "

$ws.Range("E46").Value = "  -2.07%  "

$ws.Range("E47").Value = "  -5.42%  "

$ws.Range("D48").Value = "'7.08"
$ws.Range("E48").Value = "  -3.55%  "

$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("D50").Value = "'2.228.51"
$ws.Range("E50").Value = "  -2.32%  "

$ws.Range("D51").Value = "'3.56"
$ws.Range("E51").Value = "  -15.16%  "
